$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2711.151890339096
$ws.Range("C2").Value = 12147.56341574007
$ws.Range("D2").Value = 16294.79002210704
